# Update the "as_of_utc" timestamp column (AA) for every data row on the
# "Главные" and "Линейные" sheets from 2025-12-11 03:02:38 to 2025-12-11 07:02:19.

$wb = $excel.ActiveWorkbook

$oldStamp = "2025-12-11 03:02:38"
$newStamp = "2025-12-11 07:02:19"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($row = 2; $row -le 26; $row++) {
        $cell = $ws.Cells.Item($row, 27)  # column AA = 27
        if ($cell.Value2 -eq $oldStamp) {
            $cell.Value = $newStamp
        }
    }
}
